$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 26) documenting the Labarrera et al. (2013) paper
$ws.Range("A26").Value = "Labarrera et al. (2013)"
$ws.Range("B26").Value = "Tonic inhibition sets the state of excitability in olfactorybulb granule cells"
$ws.Range("C26").Value = "Mice"
$ws.Range("D26").Value = "GC"
$ws.Range("E26").Value = "Tonic, bulb wide inhibition"
$ws.Range("I26").Value = "STD"
$ws.Range("J26").Value = 37.5
$ws.Range("K26").Value = "Not reported"

# Match the author's final cursor position/selection
$ws.Range("L26").Select()
